$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) ---
$ws.Range("A1").Value2 = "mx_state"
$ws.Range("B1").Value2 = "mx_municipality"
$ws.Range("C1").Value2 = "n_matriculas"
$ws.Range("D1").Value2 = "pct_matriculas"

# --- 2. Convert UPPERCASE state/municipality/label text (columns A & B,
#        rows 2-75) to Title Case ---
$ti = (Get-Culture).TextInfo

for ($r = 2; $r -le 75; $r++) {
    foreach ($col in @(1, 2)) {
        $cell = $ws.Cells.Item($r, $col)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val -ne "") {
            $cell.Value2 = $ti.ToTitleCase($val.ToLower())
        }
    }
}

# --- 3. Remove trailing metadata/footer rows (77-81) ---
$ws.Rows("77:81").Delete()
